$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / rId1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5552
$ws1.Range("F8").Value = 910
$ws1.Range("F9").Value = 145
$ws1.Range("F10").Value = 2483
$ws1.Range("F12").Value = 104
$ws1.Range("F14").Value = 71
$ws1.Range("F15").Value = 6
$ws1.Range("F16").Value = 2324
$ws1.Range("F17").Value = 270

# Sheet "全部类型" (sheet4 / rId4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5552
$ws4.Range("F10").Value = 910
$ws4.Range("F11").Value = 145
$ws4.Range("F12").Value = 2483
$ws4.Range("F14").Value = 104
$ws4.Range("F17").Value = 71
$ws4.Range("F18").Value = 6
$ws4.Range("F19").Value = 2324
$ws4.Range("F20").Value = 270
